# 1. Rename header cells in existing sheets
$wb = $excel.ActiveWorkbook
$wsWeekly = $wb.Worksheets.Item(1)
$wsWeekly.Range("B1").Value = "Weekly_PO_Qty"

$wsMonthly = $wb.Worksheets.Item(2)
$wsMonthly.Range("B1").Value = "Monthly_PO_Qty"

# 2. Add the new "PO Forecast" worksheet at the end of the workbook
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$newSheet = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $lastSheet)
$newSheet.Name = "PO Forecast"

# Header row
$newSheet.Range("A1").Value = "ds"
$newSheet.Range("B1").Value = "PO_Forecast"
$newSheet.Range("C1").Value = "yhat_lower"
$newSheet.Range("D1").Value = "yhat_upper"

# Match the header formatting used on the other sheets (bold/border/centered)
$wsWeekly.Range("A1:B1").Copy()
$newSheet.Range("A1:D1").PasteSpecial(-4122)

# Data rows (ds, PO_Forecast, yhat_lower, yhat_upper)
$data = @(
    @(45137.99999999999,189,-306.0910904943773,723.994414250933),
    @(45151.99999999999,192,-336.8435437455274,761.9769652293602),
    @(45158.99999999999,193,-321.9191877549059,733.3661865738694),
    @(45165.99999999999,195,-328.6178800563654,749.1341967663283),
    @(45172.99999999999,196,-380.7522533986624,715.3456818996628),
    @(45179.99999999999,197,-335.2636370263596,730.3867010146049),
    @(45186.99999999999,199,-328.6779349216671,706.198807283359),
    @(45193.99999999999,200,-340.8259390026039,699.3179772222371),
    @(45200.99999999999,202,-349.7843597711542,752.8180000419035),
    @(45207.99999999999,203,-318.7622797404824,753.7657176699353),
    @(45214.99999999999,205,-312.0292892113168,756.9317039157346),
    @(45221.99999999999,206,-319.2041769938524,744.0481983129603),
    @(45228.99999999999,208,-290.2998630578529,732.4602209182495),
    @(45235.99999999999,209,-331.8705021760244,735.8090580796667),
    @(45242.99999999999,211,-340.2847335446766,681.987128658388),
    @(45249.99999999999,212,-350.0279047347225,736.6011150360735),
    @(45256.99999999999,214,-335.8592905336717,790.0278316330478),
    @(45263.99999999999,215,-309.3602499881044,734.1610556654952),
    @(45277.99999999999,218,-313.9091951875596,774.0068223143919),
    @(45298.99999999999,222,-342.5987351684619,762.3776972929643),
    @(45312.99999999999,225,-297.7005089952436,754.7290995676346),
    @(45319.99999999999,227,-327.0061353421937,732.8874421805135),
    @(45368.99999999999,237,-267.0890865424281,756.3994972889365),
    @(45375.99999999999,239,-259.6197744385266,737.3762035759759),
    @(45403.99999999999,245,-295.617287479181,757.5162455212508),
    @(45410.99999999999,246,-295.0955570585316,781.5017762919539),
    @(45417.99999999999,247,-295.1806868613756,737.8532187224178),
    @(45424.99999999999,249,-278.3773641172883,773.6906228974606),
    @(45431.99999999999,250,-313.4854370785845,796.055673652706),
    @(45438.99999999999,252,-245.2348685276074,798.2701754189258),
    @(45445.99999999999,253,-267.9118403851724,738.9139553847103),
    @(45452.99999999999,255,-300.904612020534,788.8711030347267),
    @(45459.99999999999,256,-299.483195415236,804.7090664981512),
    @(45473.99999999999,259,-331.855199547621,811.0586855005876),
    @(45487.99999999999,262,-255.3045604640578,814.6153599073094),
    @(45494.99999999999,264,-285.5728322588703,806.7286987307905),
    @(45585.99999999999,283,-260.3448020535126,800.0379097409231),
    @(45592.99999999999,284,-249.3646944803042,849.4704125435342),
    @(45599.99999999999,286,-249.3421762474821,818.6008478827405),
    @(45606.99999999999,287,-200.0433017470792,828.3835339213023),
    @(45634.99999999999,293,-240.1813580117824,842.3820942152601),
    @(45641.99999999999,295,-229.8953842520768,785.9992351505401),
    @(45648.99999999999,296,-231.2311762858348,845.9187546527149),
    @(45655.99999999999,297,-228.3951548927876,830.3495367545362),
    @(45662.99999999999,299,-200.9840909515363,847.4524696223967),
    @(45669.99999999999,300,-301.2135776293399,826.1256406448408),
    @(45676.99999999999,302,-242.3401270123555,872.9022327239651),
    @(45683.99999999999,303,-235.4991382778352,797.3157051799433),
    @(45690.99999999999,305,-253.87187463073,827.3435958028618),
    @(45697.99999999999,306,-206.6328201773559,849.3587699982918)
)

$r = 2
foreach ($row in $data) {
    $newSheet.Cells.Item($r, 1).Value = $row[0]
    $newSheet.Cells.Item($r, 2).Value = $row[1]
    $newSheet.Cells.Item($r, 3).Value = $row[2]
    $newSheet.Cells.Item($r, 4).Value = $row[3]
    $r = $r + 1
}

# Match the date-column formatting (column A) used on the other sheets
$wsWeekly.Range("A2").Copy()
$newSheet.Range("A2:A51").PasteSpecial(-4122)

# Match the page margins used on the other sheets (in points: 1in = 72pt)
$newSheet.PageSetup.LeftMargin = 54
$newSheet.PageSetup.RightMargin = 54
$newSheet.PageSetup.TopMargin = 72
$newSheet.PageSetup.BottomMargin = 72
$newSheet.PageSetup.HeaderMargin = 36
$newSheet.PageSetup.FooterMargin = 36

# Keep the originally active sheet selected, as before this edit
$wsWeekly.Activate()
